$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.965.80"
$ws.Range("E2").Value = "  +6.53%  "
$ws.Range("D3").Value = "3.018.03"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +12.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.014.18"
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.73%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("E11").Value = "  +6.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.73%  "
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.39%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "65.915.86"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").Value = "3.515.76"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("E18").Value = "  +7.67%  "
$ws.Range("D19").Value = "3.011.85"
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.66%  "
$ws.Range("E21").Value = "  +8.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.690"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.84%  "
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("E25").Value = "  +12.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.07"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +15.84%  "
$ws.Range("E30").Value = "  +15.96%  "
$ws.Range("E31").Value = "  -6.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.85"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.82%  "
$ws.Range("E38").Value = "  +14.56%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("E41").Value = "  +16.48%  "
$ws.Range("E42").Value = "  +6.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.39"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.80%  "
$ws.Range("D46").Value = "2.806.96"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("E47").Value = "  +5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.82"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.80%  "
$ws.Range("E51").Value = "  +4.57%  "
